# Refresh the cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# values on the active worksheet to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.091.42"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "1.667.36"
$ws.Range("E3").Value = "  -1.27%  "

$ws.Range("D5").Value = "'209.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.68%  "

$ws.Range("D6").Value = "'0.5245"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.01%  "

$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.42%  "

$ws.Range("D8").Value = "'0.2616"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.70%  "

$ws.Range("D9").Value = "'0.06288"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.95%  "

$ws.Range("D10").Value = "'21.08"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.73%  "

$ws.Range("D11").Value = "'0.07528"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.84%  "

$ws.Range("D12").Value = "1.665.34"
$ws.Range("E12").Value = "  -1.55%  "

$ws.Range("D13").Value = "'4.426"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.92%  "

$ws.Range("D14").Value = "'0.5510"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.62%  "

$ws.Range("D15").Value = "'66.35"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.13%  "

$ws.Range("D16").Value = "'0.000007924"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.91%  "

$ws.Range("D17").Value = "26.127.45"
$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").Value = "'1.003"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.51%  "

$ws.Range("D19").Value = "'4.707"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.67%  "

$ws.Range("D20").Value = "'186.08"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.11%  "

$ws.Range("D21").Value = "'10.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.46%  "

$ws.Range("D22").Value = "'6.159"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.25%  "

$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.46%  "

$ws.Range("D24").Value = "'149.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("D25").Value = "'0.1239"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.45%  "

$ws.Range("D26").Value = "'7.447"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.82%  "

$ws.Range("D27").Value = "'15.88"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.85%  "

$ws.Range("D28").Value = "'0.06343"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.43%  "

$ws.Range("D29").Value = "'1.352"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.88%  "

$ws.Range("D30").Value = "'1.275"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.59%  "

$ws.Range("D31").Value = "'3.487"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.62%  "

$ws.Range("D32").Value = "'3.408"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.76%  "

$ws.Range("D33").Value = "'1.632"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.97%  "

$ws.Range("D34").Value = "'0.9991"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.06%  "

$ws.Range("D35").Value = "'2.409"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("D36").Value = "'0.6012"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.67%  "

$ws.Range("D37").Value = "'2.731"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.06%  "

$ws.Range("D38").Value = "1.109.90"
$ws.Range("E38").Value = "  +0.45%  "

$ws.Range("D39").Value = "'6.089"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("D40").Value = "'0.01614"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.62%  "

$ws.Range("D41").Value = "'0.8699"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.93%  "

$ws.Range("D43").Value = "'99.77"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("D44").Value = "1.817.67"
$ws.Range("E44").Value = "  -1.17%  "

$ws.Range("D45").Value = "'0.00000000106"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.01%  "

$ws.Range("D46").Value = "'55.23"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.07%  "

$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").Value = "'8.012"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("D49").Value = "'0.05231"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("D50").Value = "'0.4244"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("D51").Value = "'5.924"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.06%  "
